# "added a few things, and glow texture instructions"
#
# After the existing "Might want to remake the stair system..." paragraph,
# add two new paragraphs (the "bumper" note and the glow-effect note),
# separated by a blank paragraph. The hidden "_GoBack" bookmark that Word
# keeps at the end of the document must end up wrapping the very end of the
# new final paragraph (exactly like it originally wrapped the end of the
# last paragraph before the edit).

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document/holds the bookmark.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Remove the existing hidden "_GoBack" bookmark up front - we'll recreate a
# fresh one in the right spot once the new trailing paragraph exists. (Word
# auto-manages "_GoBack"; re-adding it by name after the edit mirrors the
# normal behaviour of the bookmark following the most recent edit point.)
$goBack = $lastRange.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Add the two new paragraphs (plus the blank separator) right after it.
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()
$tail.InsertParagraphAfter()
$tail.InsertParagraphAfter()

$bumperPara = $d.Paragraphs.Item($lastPara.Index + 1)
$bumperPara.Range.Text = 'All walls should have a "bumper" section on them when they reach the ceiling, floor, or a corner'

$glowPara = $d.Paragraphs.Item($lastPara.Index + 3)
# Temporarily tack on a unique marker so the final insertion point is not
# the literal last character of the whole document (that special position
# makes a freshly-added collapsed bookmark expand to cover the paragraph
# instead of staying collapsed). We'll strip the marker back out below.
$marker = "ZZZ_GOBACK_MARKER_ZZZ"
$glowPara.Range.Text = "For the glow effect, add the emissions channel in substance painter and whatever you paint with that will be the glow effect. " + $marker

# Re-create "_GoBack" collapsed right before the marker (i.e. at the true
# end of the visible glow-effect text).
$markerRange = $d.Content
$markerRange.Find.Execute($marker, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
$markerRange.Collapse(1)
$markerRange.Bookmarks.Add("_GoBack")

# Now remove the temporary marker text itself.
$cleanup = $d.Content
$cleanup.Find.Execute($marker, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
$cleanup.Text = ""
